# Update the 20x5 division-problem table: each populated cell's text is
# replaced by addressing the cell directly via Table.Cell(row, col), so the
# substitution is unambiguous even where an old/new value reappears
# elsewhere in the table (e.g. "76÷8=" and "19÷6=" show up as both a
# source and a target value in different cells).
$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Row 1 (table rows 1,5,9,13,17 contain data; the others are blank spacer rows)
$t.Cell(1,1).Range.Text = "11÷2="
$t.Cell(1,2).Range.Text = "99÷9="
$t.Cell(1,3).Range.Text = "47÷4="
$t.Cell(1,4).Range.Text = "54÷3="
$t.Cell(1,5).Range.Text = "51÷3="

# Row 5
$t.Cell(5,1).Range.Text = "15÷2="
$t.Cell(5,2).Range.Text = "83÷5="
$t.Cell(5,3).Range.Text = "19÷6="
$t.Cell(5,4).Range.Text = "24÷3="
$t.Cell(5,5).Range.Text = "13÷2="

# Row 9
$t.Cell(9,1).Range.Text = "85÷7="
$t.Cell(9,2).Range.Text = "76÷8="
$t.Cell(9,3).Range.Text = "48÷2="
$t.Cell(9,4).Range.Text = "28÷3="
$t.Cell(9,5).Range.Text = "15÷7="

# Row 13
$t.Cell(13,1).Range.Text = "56÷4="
$t.Cell(13,2).Range.Text = "30÷6="
$t.Cell(13,3).Range.Text = "40÷3="
$t.Cell(13,4).Range.Text = "50÷9="
$t.Cell(13,5).Range.Text = "23÷8="

# Row 17
$t.Cell(17,1).Range.Text = "31÷5="
$t.Cell(17,2).Range.Text = "97÷8="
$t.Cell(17,3).Range.Text = "16÷3="
$t.Cell(17,4).Range.Text = "56÷9="
$t.Cell(17,5).Range.Text = "43÷3="
